$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H29 and H30: mark as corrected in version 3.35
$ws.Range("H29").Value = "3.35"
$ws.Range("H30").Value = "3.35"

# Row 31: new entry - error with servo addressing
$ws.Range("A31").Value = 1040
$ws.Range("B31").Value2 = 43955
$ws.Range("C31").Value = "Fehler bei Adressirung von Servos behoben"
$ws.Range("D31").Value = "Harold"
$ws.Range("E31").Value = "Fehler"
$ws.Range("H31").Value = "3.36"

# Row 32: new entry - automatic comport list update
$ws.Range("A32").Value = 1041
$ws.Range("B32").Value2 = 43955
$ws.Range("C32").Value = "Automatisches Update der Comport-Liste, wenn der Einstellungstab geöffnet wird"
$ws.Range("D32").Value = "Harold"
$ws.Range("E32").Value = "Neue Funktion"
$ws.Range("H32").Value = "3.36"
$ws.Rows.Item(32).RowHeight = 30

# B31/B32 picked up the default numeric/general style from the empty template
# cells; restore the date display format used by the rest of column B.
$ws.Range("B30").Copy()
$ws.Range("B31:B32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update selection to reflect the edited workbook's active cell
$ws.Range("H32").Select()
